$d = $word.ActiveDocument

function Replace-Unique($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

# Header date line
Replace-Unique "2026-01-07 Wednesday" "2026-01-08 Thursday"

# Table data - row 1 (unique strings, safe to use document-wide replace)
Replace-Unique "22÷5=4, 2" "88÷7=12, 4"
Replace-Unique "65÷5=13, 0" "15÷7=2, 1"
Replace-Unique "64÷5=12, 4" "27÷5=5, 2"
Replace-Unique "73÷5=14, 3" "25÷9=2, 7"
Replace-Unique "12÷4=3, 0" "54÷6=9, 0"

# Table data - row 5 (unique strings)
Replace-Unique "68÷3=22, 2" "35÷8=4, 3"
Replace-Unique "15÷8=1, 7" "40÷7=5, 5"
Replace-Unique "95÷2=47, 1" "98÷5=19, 3"
Replace-Unique "49÷7=7, 0" "59÷2=29, 1"
Replace-Unique "66÷5=13, 1" "30÷4=7, 2"

# Table data - row 9 (unique strings, except col 5)
Replace-Unique "23÷9=2, 5" "89÷6=14, 5"
Replace-Unique "56÷3=18, 2" "79÷6=13, 1"
Replace-Unique "23÷6=3, 5" "30÷4=7, 2"
Replace-Unique "24÷6=4, 0" "49÷3=16, 1"

# Table data - row 13 (unique strings, except col 1)
Replace-Unique "39÷4=9, 3" "12÷2=6, 0"
Replace-Unique "40÷5=8, 0" "37÷9=4, 1"
Replace-Unique "53÷7=7, 4" "42÷7=6, 0"
Replace-Unique "46÷4=11, 2" "59÷5=11, 4"

# Table data - row 17 (unique strings)
Replace-Unique "60÷8=7, 4" "47÷9=5, 2"
Replace-Unique "91÷6=15, 1" "24÷5=4, 4"
Replace-Unique "17÷2=8, 1" "13÷7=1, 6"
Replace-Unique "92÷9=10, 2" "58÷5=11, 3"
Replace-Unique "13÷4=3, 1" "62÷4=15, 2"

# "78÷6=13, 0" appears twice (row 9 col 5, row 13 col 1) with different
# replacement targets, so scope the Find/Replace to each individual table
# cell's Range and use wdReplaceOne (not wdReplaceAll) so the match is
# confined to that cell instead of leaking across the whole document.
$t = $d.Tables.Item(1)

$cellR9C5 = $t.Cell(9, 5)
$cellR9C5.Range.Find.Execute("78÷6=13, 0", $true, $false, $false, $false, $false, $true, 0, $false, "21÷6=3, 3", 1) | Out-Null

$cellR13C1 = $t.Cell(13, 1)
$cellR13C1.Range.Find.Execute("78÷6=13, 0", $true, $false, $false, $false, $false, $true, 0, $false, "11÷8=1, 3", 1) | Out-Null
